$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Bayesian A/B Testing with Expected Loss"
$ws.Range("E3").Value = "https://lumiamitie.github.io/data/bayesian-ab-testing/"

$ws.Range("D28").Value = "Continual Learning for Robotics(2)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/143"

$ws.Range("D36").Value = "Deep Learning for Tabular Dataset"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/327"
